$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.052410633050744
$ws.Range("D2").Value = 1.05882101653686
$ws.Range("E2").Value = 1.059444172564484
$ws.Range("F2").Value = 1.070699348676744
$ws.Range("I2").Value = 1.046354416577647
$ws.Range("J2").Value = 1.057433380636787
$ws.Range("K2").Value = 1.061552204676095
$ws.Range("L2").Value = 1.062173660132338
$ws.Range("M2").Value = 1.073398489800253
$ws.Range("N2").Value = 1.022861484814707

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.05363592763987
$ws.Range("D3").Value = 1.059790069819249
$ws.Range("E3").Value = 1.060510131537757
$ws.Range("F3").Value = 1.071794264610962
$ws.Range("I3").Value = 1.046668289042889
$ws.Range("J3").Value = 1.058307806429355
$ws.Range("K3").Value = 1.06233497373276
$ws.Range("L3").Value = 1.063053212892584
$ws.Range("M3").Value = 1.074309128906949
$ws.Range("N3").Value = 1.023159309003675

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.054428602001437
$ws.Range("D4").Value = 1.06041668855398
$ws.Range("E4").Value = 1.061199976530561
$ws.Range("F4").Value = 1.072502710990073
$ws.Range("I4").Value = 1.046869603793193
$ws.Range("J4").Value = 1.058872918845883
$ws.Range("K4").Value = 1.062840434686546
$ws.Range("L4").Value = 1.063621840152145
$ws.Range("M4").Value = 1.074897738777412
$ws.Range("N4").Value = 1.023351609821489

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.054761802604691
$ws.Range("D5").Value = 1.060680018271911
$ws.Range("E5").Value = 1.061490011677388
$ws.Range("F5").Value = 1.072800533633042
$ws.Range("I5").Value = 1.046953810212442
$ws.Range("J5").Value = 1.059110325519867
$ws.Range("K5").Value = 1.063052681419723
$ws.Range("L5").Value = 1.06386077159235
$ws.Range("M5").Value = 1.075145039161107
$ws.Range("N5").Value = 1.02343235481995

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.054817746175314
$ws.Range("D6").Value = 1.060724226596678
$ws.Range("E6").Value = 1.061538711303281
$ws.Range("F6").Value = 1.072850538905271
$ws.Range("I6").Value = 1.04696792384931
$ws.Range("J6").Value = 1.059150177402083
$ws.Range("K6").Value = 1.063088303994059
$ws.Range("L6").Value = 1.063900882249534
$ws.Range("M6").Value = 1.075186553137202
$ws.Range("N6").Value = 1.023445906504309

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.054433054397996
$ws.Range("D7").Value = 1.060420207573894
$ws.Range("E7").Value = 1.061203851897882
$ws.Range("F7").Value = 1.072506690542157
$ws.Range("I7").Value = 1.046870730638266
$ws.Range("J7").Value = 1.058876091738952
$ws.Range("K7").Value = 1.062843271715605
$ws.Range("L7").Value = 1.063625033234192
$ws.Range("M7").Value = 1.074901043810517
$ws.Range("N7").Value = 1.023352689125963

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.052824764009626
$ws.Range("D8").Value = 1.059148600688907
$ws.Range("E8").Value = 1.05980439852136
$ws.Range("F8").Value = 1.071069388849474
$ws.Range("I8").Value = 1.04646086031387
$ws.Range("J8").Value = 1.057729042306945
$ws.Range("K8").Value = 1.061816961482385
$ws.Range("L8").Value = 1.062471013549757
$ws.Range("M8").Value = 1.073706376355858
$ws.Range("N8").Value = 1.022962221194039

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.049989345443611
$ws.Range("D9").Value = 1.056904595554941
$ws.Range("E9").Value = 1.057339099476137
$ws.Range("F9").Value = 1.068536358545457
$ws.Range("I9").Value = 1.04572495906196
$ws.Range("J9").Value = 1.055702396781951
$ws.Range("K9").Value = 1.060000454646397
$ws.Range("L9").Value = 1.060433602401648
$ws.Range("M9").Value = 1.071596327920691
$ws.Range("N9").Value = 1.022271006670483

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.048098008274269
$ws.Range("D10").Value = 1.055406350946041
$ws.Range("E10").Value = 1.055695995368305
$ws.Range("F10").Value = 1.066847404595085
$ws.Range("I10").Value = 1.045225157880776
$ws.Range("J10").Value = 1.054347611736267
$ws.Range("K10").Value = 1.058784016273562
$ws.Range("L10").Value = 1.059072667589385
$ws.Range("M10").Value = 1.070186292293277
$ws.Range("N10").Value = 1.021808057693352

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.047278761434493
$ws.Range("D11").Value = 1.054757052624044
$ws.Range("E11").Value = 1.054984602038404
$ws.Range("F11").Value = 1.066115993136683
$ws.Range("I11").Value = 1.045006551689061
$ws.Range("J11").Value = 1.05376008620811
$ws.Range("K11").Value = 1.058255985630782
$ws.Range("L11").Value = 1.058482724461339
$ws.Range("M11").Value = 1.069574928345286
$ws.Range("N11").Value = 1.021607084691959

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.04697441146663
$ws.Range("D12").Value = 1.054515790905988
$ws.Range("E12").Value = 1.054720369933649
$ws.Range("F12").Value = 1.06584430066341
$ws.Range("I12").Value = 1.04492502227038
$ws.Range("J12").Value = 1.053541717355593
$ws.Range("K12").Value = 1.058059654488505
$ws.Range("L12").Value = 1.058263494600363
$ws.Range("M12").Value = 1.069347717873348
$ws.Range("N12").Value = 1.021532356943342

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.047039697651604
$ws.Range("D13").Value = 1.054567546154214
$ws.Range("E13").Value = 1.054777048136345
$ws.Range("F13").Value = 1.065902580259872
$ws.Range("I13").Value = 1.044942525518714
$ws.Range("J13").Value = 1.053588564377892
$ws.Range("K13").Value = 1.058101777118404
$ws.Range("L13").Value = 1.058310524635452
$ws.Range("M13").Value = 1.069396460860242
$ws.Range("N13").Value = 1.021548389808478

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.047253604707138
$ws.Range("D14").Value = 1.054737111562561
$ws.Range("E14").Value = 1.054962760311721
$ws.Range("F14").Value = 1.066093535238266
$ws.Range("I14").Value = 1.044999819163199
$ws.Range("J14").Value = 1.053742038546654
$ws.Range("K14").Value = 1.05823976086302
$ws.Range("L14").Value = 1.05846460487312
$ws.Range("M14").Value = 1.069556149565863
$ws.Range("N14").Value = 1.021600909254472

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.04738539384235
$ws.Range("D15").Value = 1.05484157532408
$ws.Range("E15").Value = 1.055077185140673
$ws.Range("F15").Value = 1.066211187052961
$ws.Range("I15").Value = 1.045035076012697
$ws.Range("J15").Value = 1.053836581025234
$ws.Range("K15").Value = 1.05832475106033
$ws.Range("L15").Value = 1.058559525677361
$ws.Range("M15").Value = 1.069654522758232
$ws.Range("N15").Value = 1.021633257948128

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.04815237277685
$ws.Range("D16").Value = 1.055449431095212
$ws.Range("E16").Value = 1.055743209816993
$ws.Range("F16").Value = 1.066895944109834
$ws.Range("I16").Value = 1.045239619893901
$ws.Range("J16").Value = 1.054386584910358
$ws.Range("K16").Value = 1.058819032350797
$ws.Range("L16").Value = 1.059111806391858
$ws.Range("M16").Value = 1.070226849361605
$ws.Range("N16").Value = 1.02182138477065

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.048633400346967
$ws.Range("D17").Value = 1.055830575401252
$ws.Range("E17").Value = 1.056161010495177
$ws.Range("F17").Value = 1.067325450997629
$ws.Range("I17").Value = 1.045367338479081
$ws.Range("J17").Value = 1.054731347374353
$ws.Range("K17").Value = 1.059128731895526
$ws.Range("L17").Value = 1.059458062823221
$ws.Range("M17").Value = 1.070585637290569
$ws.Range("N17").Value = 1.021939254154816

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.048913948216473
$ws.Range("D18").Value = 1.056052837571767
$ws.Range("E18").Value = 1.056404714517583
$ws.Range("F18").Value = 1.067575967147656
$ws.Range("I18").Value = 1.045441623396969
$ws.Range("J18").Value = 1.054932355253793
$ws.Range("K18").Value = 1.059309248576308
$ws.Range("L18").Value = 1.059659965765732
$ws.Range("M18").Value = 1.070794834217871
$ws.Range("N18").Value = 1.022007955895657

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.049009603217148
$ws.Range("D19").Value = 1.056128614292766
$ws.Range("E19").Value = 1.056487812676602
$ws.Range("F19").Value = 1.067661385367451
$ws.Range("I19").Value = 1.045466916829918
$ws.Range("J19").Value = 1.055000879158348
$ws.Range("K19").Value = 1.059370778787871
$ws.Range("L19").Value = 1.059728798871984
$ws.Range("M19").Value = 1.070866151766407
$ws.Range("N19").Value = 1.022031373032131

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.048581793490623
$ws.Range("D20").Value = 1.055789687699493
$ws.Range("E20").Value = 1.056116183623234
$ws.Range("F20").Value = 1.06727936980533
$ws.Range("I20").Value = 1.04535365733001
$ws.Range("J20").Value = 1.054694366554506
$ws.Range("K20").Value = 1.059095517088088
$ws.Range("L20").Value = 1.059420919266777
$ws.Range("M20").Value = 1.070547150818258
$ws.Range("N20").Value = 1.021926613005843

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.047190615677759
$ws.Range("D21").Value = 1.05468718108676
$ws.Range("E21").Value = 1.054908072409431
$ws.Range("F21").Value = 1.066037304144056
$ws.Range("I21").Value = 1.044982956708299
$ws.Range("J21").Value = 1.053696847985967
$ws.Range("K21").Value = 1.05819913353718
$ws.Range("L21").Value = 1.058419234817826
$ws.Range("M21").Value = 1.069509128617692
$ws.Range("N21").Value = 1.021585445726833

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.046315664209475
$ws.Range("D22").Value = 1.053993508326826
$ws.Range("E22").Value = 1.054148548285714
$ws.Range("F22").Value = 1.065256288358661
$ws.Range("I22").Value = 1.044747976783444
$ws.Range("J22").Value = 1.053068882636585
$ws.Range("K22").Value = 1.057634401035225
$ws.Range("L22").Value = 1.057788864331591
$ws.Range("M22").Value = 1.068855772516386
$ws.Range("N22").Value = 1.021370492171343

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.046779518228064
$ws.Range("D23").Value = 1.05436128351041
$ws.Range("E23").Value = 1.054551180794331
$ws.Range("F23").Value = 1.065670327566077
$ws.Range("I23").Value = 1.044872724888932
$ws.Range("J23").Value = 1.053401853873219
$ws.Range("K23").Value = 1.057933884876455
$ws.Range("L23").Value = 1.058123090232841
$ws.Range("M23").Value = 1.069202196651117
$ws.Range("N23").Value = 1.021484485714395

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.048605112492233
$ws.Range("D24").Value = 1.055808163255376
$ws.Range("E24").Value = 1.056136438932662
$ws.Range("F24").Value = 1.067300191937296
$ws.Range("I24").Value = 1.045359839904922
$ws.Range("J24").Value = 1.054711076861702
$ws.Range("K24").Value = 1.059110525818727
$ws.Range("L24").Value = 1.059437703034554
$ws.Range("M24").Value = 1.070564541439589
$ws.Range("N24").Value = 1.021932325149851

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.050722547293474
$ws.Range("D25").Value = 1.057485116605341
$ws.Range("E25").Value = 1.057976358906149
$ws.Range("F25").Value = 1.069191250888794
$ws.Range("I25").Value = 1.04591682746964
$ws.Range("J25").Value = 1.056226978936944
$ws.Range("K25").Value = 1.060471019730468
$ws.Range("L25").Value = 1.060960787204833
$ws.Range("M25").Value = 1.072142410604031
$ws.Range("N25").Value = 1.022450078181999
